$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing Chinese full-stop "。" from the "Des" (描述) column
# for each level entry (D3:D7).
$ws.Range("D3").Value = "新手试炼"
$ws.Range("D4").Value = "哥布林的栖息地"
$ws.Range("D5").Value = "可怕的关卡，充满危机"
$ws.Range("D6").Value = "复杂的深林住着远古的精灵"
$ws.Range("D7").Value = "充满绝望的气息"

# D4 was previously using the default format; bring it in line with the
# rest of the "Des" column (D3, D5:D7) by copying that cell format over.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active cell selection to D7 (matches the saved selection state).
$ws.Range("D7").Select()
